# Natmi following Dr Hou advice
#
# Expands the Cck -> Cckbr LR-pair sheet from two sending clusters
# (FAPs, sCs -> FAPs) to four sending clusters (ECs, FAPs, M2, sCs),
# all signalling Cck -> Cckbr onto the FAPs target cluster, and
# refreshes every numeric (expression / specificity) column to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> FAPs (Cck -> Cckbr)
$ws.Cells.Item(2, 1).Value  = "ECs"
$ws.Cells.Item(2, 2).Value  = "Cck"
$ws.Cells.Item(2, 3).Value  = "Cckbr"
$ws.Cells.Item(2, 4).Value  = "FAPs"
$ws.Cells.Item(2, 5).Value  = 1
$ws.Cells.Item(2, 6).Value  = 0.3333333333333333
$ws.Cells.Item(2, 7).Value  = 0.1295
$ws.Cells.Item(2, 8).Value  = 0.3885
$ws.Cells.Item(2, 9).Value  = 0.03378978388280364
$ws.Cells.Item(2, 10).Value = 0.03378978388280364
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.4118343333333334
$ws.Cells.Item(2, 14).Value = 1.235503
$ws.Cells.Item(2, 15).Value = 1
$ws.Cells.Item(2, 16).Value = 1
$ws.Cells.Item(2, 17).Value = 0.05333254616666667
$ws.Cells.Item(2, 18).Value = 0.4799929155
$ws.Cells.Item(2, 19).Value = 0.03378978388280364
$ws.Cells.Item(2, 20).Value = 0.03378978388280364

# Row 3: FAPs -> FAPs (Cck -> Cckbr)
$ws.Cells.Item(3, 1).Value  = "FAPs"
$ws.Cells.Item(3, 2).Value  = "Cck"
$ws.Cells.Item(3, 3).Value  = "Cckbr"
$ws.Cells.Item(3, 4).Value  = "FAPs"
$ws.Cells.Item(3, 5).Value  = 3
$ws.Cells.Item(3, 6).Value  = 1
$ws.Cells.Item(3, 7).Value  = 1.300224333333333
$ws.Cells.Item(3, 8).Value  = 3.900673
$ws.Cells.Item(3, 9).Value  = 0.3392609978571102
$ws.Cells.Item(3, 10).Value = 0.3392609978571102
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.4118343333333334
$ws.Cells.Item(3, 14).Value = 1.235503
$ws.Cells.Item(3, 15).Value = 1
$ws.Cells.Item(3, 16).Value = 1
$ws.Cells.Item(3, 17).Value = 0.5354770215021111
$ws.Cells.Item(3, 18).Value = 4.819293193519
$ws.Cells.Item(3, 19).Value = 0.3392609978571102
$ws.Cells.Item(3, 20).Value = 0.3392609978571102

# Row 4: M2 -> FAPs (Cck -> Cckbr)
$ws.Cells.Item(4, 1).Value  = "M2"
$ws.Cells.Item(4, 2).Value  = "Cck"
$ws.Cells.Item(4, 3).Value  = "Cckbr"
$ws.Cells.Item(4, 4).Value  = "FAPs"
$ws.Cells.Item(4, 5).Value  = 1
$ws.Cells.Item(4, 6).Value  = 0.3333333333333333
$ws.Cells.Item(4, 7).Value  = 0.1462553333333333
$ws.Cells.Item(4, 8).Value  = 0.438766
$ws.Cells.Item(4, 9).Value  = 0.03816166876479336
$ws.Cells.Item(4, 10).Value = 0.03816166876479336
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.4118343333333334
$ws.Cells.Item(4, 14).Value = 1.235503
$ws.Cells.Item(4, 15).Value = 1
$ws.Cells.Item(4, 16).Value = 1
$ws.Cells.Item(4, 17).Value = 0.06023296769977778
$ws.Cells.Item(4, 18).Value = 0.542096709298
$ws.Cells.Item(4, 19).Value = 0.03816166876479336
$ws.Cells.Item(4, 20).Value = 0.03816166876479336

# Row 5: sCs -> FAPs (Cck -> Cckbr)
$ws.Cells.Item(5, 1).Value  = "sCs"
$ws.Cells.Item(5, 2).Value  = "Cck"
$ws.Cells.Item(5, 3).Value  = "Cckbr"
$ws.Cells.Item(5, 4).Value  = "FAPs"
$ws.Cells.Item(5, 5).Value  = 3
$ws.Cells.Item(5, 6).Value  = 1
$ws.Cells.Item(5, 7).Value  = 2.256539666666666
$ws.Cells.Item(5, 8).Value  = 6.769619
$ws.Cells.Item(5, 9).Value  = 0.5887875494952928
$ws.Cells.Item(5, 10).Value = 0.5887875494952928
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.4118343333333334
$ws.Cells.Item(5, 14).Value = 1.235503
$ws.Cells.Item(5, 15).Value = 1
$ws.Cells.Item(5, 16).Value = 1
$ws.Cells.Item(5, 17).Value = 0.9293205092618888
$ws.Cells.Item(5, 18).Value = 8.363884583356999
$ws.Cells.Item(5, 19).Value = 0.5887875494952928
$ws.Cells.Item(5, 20).Value = 0.5887875494952928
